$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.040.57"
$ws.Range("E2").Value = "  -9.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.845.88"
$ws.Range("E3").Value = "  -9.39%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.31"
$ws.Range("E5").Value = "  -9.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "119.60"
$ws.Range("E6").Value = "  -10.78%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.843.65"
$ws.Range("E8").Value = "  -9.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -5.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  -13.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.70"
$ws.Range("E11").Value = "  -11.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.420"
$ws.Range("E12").Value = "  -6.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000208"
$ws.Range("E13").Value = "  -13.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.72"
$ws.Range("E14").Value = "  -10.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.117"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.303.92"
$ws.Range("E16").Value = "  -9.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.851.15"
$ws.Range("E17").Value = "  -9.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "57.013.02"
$ws.Range("E18").Value = "  -9.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  -5.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "404.23"
$ws.Range("E20").Value = "  -11.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("E21").Value = "  -9.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.634"
$ws.Range("E22").Value = "  -7.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.59"
$ws.Range("E23").Value = "  -12.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.25"
$ws.Range("E24").Value = "  -6.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "75.20"
$ws.Range("E25").Value = "  -8.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.39"
$ws.Range("E28").Value = "  -10.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.87"
$ws.Range("E29").Value = "  -9.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -8.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.21"
$ws.Range("E31").Value = "  -9.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.77"
$ws.Range("E32").Value = "  -12.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0907"
$ws.Range("E33").Value = "  -9.12%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.12"
$ws.Range("E34").Value = "  -5.36%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  -10.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.871"
$ws.Range("E36").Value = "  -13.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("E37").Value = "  -18.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.14"
$ws.Range("E38").Value = "  +0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0603"
$ws.Range("E39").Value = "  -16.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0332"
$ws.Range("E40").Value = "  -13.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.102"
$ws.Range("E41").Value = "  -8.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.569.49"
$ws.Range("E42").Value = "  -7.07%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "343.84"
$ws.Range("E44").Value = "  -10.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  -11.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "116.48"
$ws.Range("E46").Value = "  -7.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.222"
$ws.Range("E47").Value = "  -10.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.103"
$ws.Range("E48").Value = "  -6.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.88"
$ws.Range("E49").Value = "  -10.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.11"
$ws.Range("E50").Value = "  -10.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.88"
$ws.Range("E51").Value = "  -12.03%  "
